$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 90

# Reuse the date-formatted style already used by column A (avoids creating a
# duplicate style entry) by copying the format from the row above.
$ws.Range("A89").Copy()
$ws.Range("A90").PasteSpecial(-4122)

$ws.Cells.Item($r, 1).Value = (Get-Date -Year 2024 -Month 5 -Day 1).Date
$ws.Cells.Item($r, 2).Value = 105.07416888738
$ws.Cells.Item($r, 3).Value = 120.001606396194

# Columns D-G hold values that look numeric but must stay text, matching the
# shared-string cells used throughout the rest of the sheet.
$ws.Cells.Item($r, 4).NumberFormat = "@"
$ws.Cells.Item($r, 4).Value = "110.3"

$ws.Cells.Item($r, 5).NumberFormat = "@"
$ws.Cells.Item($r, 5).Value = "111.5"

$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = " 89.1"

$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = "164.0"
